$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PrecedenceRules")

# The three blocks of rows (43/44), (45/46), (47/48) each contain a
# "keyword rule" row (col A populated) followed by a "reciprocal" row
# (col A blank). The edit swaps the two rows within each pair, moving
# the keyword text down one row while keeping the D/F numbers attached
# to the same logical entry.

$keyword = "farming | farming$|^ferme | ferme | ferme$|farms|^farm | farm |farm$|farm,|acre|nurser|cattle| ranch|stable| sod | sod$|livestock|vineyard"

# Pair 1: rows 43 <-> 44
$ws.Range("A43").ClearContents()
$ws.Range("D43").Value = 44
$ws.Range("F43").Value = 11

$ws.Range("A44").Value = $keyword
$ws.Range("D44").Value = 11
$ws.Range("F44").Value = 44

# Pair 2: rows 45 <-> 46
$ws.Range("A45").ClearContents()
$ws.Range("D45").Value = 45
$ws.Range("F45").Value = 11

$ws.Range("A46").Value = $keyword
$ws.Range("D46").Value = 11
$ws.Range("F46").Value = 45

# Pair 3: rows 47 <-> 48
$ws.Range("A47").ClearContents()
$ws.Range("D47").Value = 81
$ws.Range("F47").Value = 11

$ws.Range("A48").Value = $keyword
$ws.Range("D48").Value = 11
$ws.Range("F48").Value = 81

# Update the view to match the saved selection/scroll position.
$ws.Activate()
$excel.Goto($ws.Range("A13"), $true)
$ws.Range("D42").Select()
